$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "22.475.57"
$ws.Range("D3").Value = "1.575.37"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  -0.15%  "
Set-TextValue $ws.Range("D6") "288.09"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  +0.99%  "
Set-TextValue $ws.Range("D8") "47.84"
$ws.Range("E8").Value = "  -2.98%  "
$ws.Range("E9").Value = "  -0.52%  "
Set-TextValue $ws.Range("D10") "1.151"
$ws.Range("E10").Value = "  +2.32%  "
Set-TextValue $ws.Range("D11") "0.07558"
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("E12").Value = "  -0.12%  "
Set-TextValue $ws.Range("D13") "20.80"
$ws.Range("E13").Value = "  +0.27%  "
Set-TextValue $ws.Range("D14") "5.949"
$ws.Range("E14").Value = "  +0.50%  "
Set-TextValue $ws.Range("D15") "6.958"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "1.567.80"
$ws.Range("E16").Value = "  +0.33%  "
Set-TextValue $ws.Range("D17") "0.00001122"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("E19").Value = "  -0.13%  "
Set-TextValue $ws.Range("D20") "0.9999"
$ws.Range("E20").Value = "  -0.17%  "
Set-TextValue $ws.Range("D21") "6.395"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("E22").Value = "  +3.18%  "
Set-TextValue $ws.Range("D23") "12.02"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").Value = "22.458.18"
$ws.Range("E24").Value = "  +0.33%  "
Set-TextValue $ws.Range("D25") "2.388"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("E26").Value = "  +4.49%  "
Set-TextValue $ws.Range("D27") "150.51"
$ws.Range("E27").Value = "  +0.63%  "
Set-TextValue $ws.Range("D28") "19.70"
$ws.Range("E28").Value = "  +0.62%  "
Set-TextValue $ws.Range("D29") "4.990"
$ws.Range("E29").Value = "  -0.15%  "
Set-TextValue $ws.Range("D30") "125.58"
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").Value = "1.747.27"
$ws.Range("E31").Value = "  +0.52%  "
Set-TextValue $ws.Range("D32") "1.091"
$ws.Range("E32").Value = "  +3.50%  "
$ws.Range("E33").Value = "  +0.54%  "
Set-TextValue $ws.Range("D34") "1.996"
$ws.Range("E34").Value = "  +0.33%  "
Set-TextValue $ws.Range("D35") "9.864"
$ws.Range("E35").Value = "  +3.47%  "
Set-TextValue $ws.Range("D36") "0.08358"
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("E37").Value = "  +3.50%  "
$ws.Range("E38").Value = "  +1.09%  "
Set-TextValue $ws.Range("D39") "0.06398"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D40") "5.375"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D41") "1.296"
$ws.Range("E41").Value = "  -0.95%  "
Set-TextValue $ws.Range("D42") "11.48"
$ws.Range("E42").Value = "  +3.23%  "
Set-TextValue $ws.Range("D43") "0.6275"
$ws.Range("E43").Value = "  +3.65%  "
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("E45").Value = "  -0.20%  "
Set-TextValue $ws.Range("D46") "0.6120"
$ws.Range("E46").Value = "  +6.87%  "
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("E48").Value = "  +2.47%  "
Set-TextValue $ws.Range("D49") "125.43"
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("E50").Value = "  +0.10%  "
Set-TextValue $ws.Range("D51") "0.07224"
$ws.Range("E51").Value = "  +0.03%  "
